$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start clean: remove all existing cell content/formatting on the sheet so we
# can rebuild the new layout (1 header row + 3 data rows, with two new
# leading columns) from scratch.
$ws.Cells.Clear()

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Header cells F1:K1 use the smaller 9pt font (same font used for text data
# cells elsewhere in the sheet).
$ws.Range("F1:K1").Font.Size = 9

# ---- Data rows (rows 2-4) ----
# Row 2: La Goule
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 510100
$ws.Range("C2").Value = "La Goule"
$ws.Range("D2").Value = 1894
$ws.Range("E2").Value = 1958
$ws.Range("F2").Value = 22
$ws.Range("G2").Value = 5.6
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 12
$ws.Range("K2").Value = 26

# Row 3: Bellerive
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 109915
$ws.Range("C3").Value = "Bellerive"
$ws.Range("D3").Value = 1905
$ws.Range("E3").Value = 2002
$ws.Range("F3").Value = 9.8
$ws.Range("G3").Value = 0.52
$ws.Range("H3").Value = 0.46
$ws.Range("I3").Value = 1.34
$ws.Range("J3").Value = 0.96
$ws.Range("K3").Value = 2.3

# Row 4: Bassecourt
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 109900
$ws.Range("C4").Value = "Bassecourt"
$ws.Range("D4").Value = 1920
$ws.Range("E4").Value = 2001
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 1.03
$ws.Range("H4").Value = 0.86
$ws.Range("I4").Value = 2.6
$ws.Range("J4").Value = 1.9
$ws.Range("K4").Value = 4.5

# All data cells use the smaller 9pt font, same as the rest of the sheet.
$ws.Range("A2:K4").Font.Size = 9

# idx / idx2 / Date Start / Date End columns are whole numbers -> "0" format.
$ws.Range("A2:B4").NumberFormat = "0"
$ws.Range("D2:E4").NumberFormat = "0"

# (m3/s), (MW1), (MW2), (GWh)* columns use two decimal places.
$ws.Range("F2:K4").NumberFormat = "0.00"

# Name column (C2:C4) keeps the default 9pt text style already applied above.

# Restore the selection to match the post-edit workbook (row 2 selected).
$ws.Range("A2:K2").Select()
